$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = 33.7285886159742
    3   = 39.03984043145672
    4   = 8.801180051587268
    5   = 48.3073919398682
    6   = 38.53995634831496
    7   = 43.85792228928054
    8   = 44.39438612824315
    9   = 0
    10  = 36.45041950800235
    11  = 22.52496801570798
    12  = 1.140786095189984
    13  = 0.326073436613155
    15  = 4.094147680108867
    93  = 0.0673803069389582
    95  = 2.084093141159541
    96  = 0.0752336908640597
    97  = 1.159196887262711
    98  = 1.031577911451526
    99  = 0.04783655377309661
    100 = 7.446331095872016
    101 = 0.5837523794174964
    102 = 2.249894246604051
    103 = 0.009259861259759665
    104 = 0.1981251808089306
    105 = 0.2053204491752422
    106 = 4.759414553207573
    107 = 1.548957220853535
    108 = 12.97948785432254
    109 = 0
    110 = 0.673890896375966
    112 = 0.01157282758218413
    113 = 3.297598604791099
    121 = 12.03227757330267
    122 = 0
    123 = 1.126081545623826
    124 = 0.5426503968743743
    125 = 178.6499685109372
    126 = 4.485537683325485
    127 = 60.08636034716793
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 3).Value = $updates[$row]
}
